# Refresh the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped values, as produced by the scheduled GitHub Actions
# update job.
#
# NOTE: Price values in column D are stored as plain text (e.g. "29.222.97",
# "1.001", "0.000007897") rather than numbers, matching how the source data
# looks (some prices even contain multiple "." thousands separators, which
# are not valid numeric literals at all). A leading apostrophe is used in the
# assigned values below to force Excel to keep them as text instead of
# silently reinterpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.222.97'
$ws.Range('D3').Value = '''1.859.61'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''0.7106'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').Value = '''237.91'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('D7').Value = '''1.001'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '''0.08125'
$ws.Range('E8').Value = '  +8.92%  '
$ws.Range('D9').Value = '''0.3042'
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').Value = '''23.17'
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('D11').Value = '''0.08187'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').Value = '''1.868.93'
$ws.Range('E12').Value = '  -6.16%  '
$ws.Range('D13').Value = '''5.172'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').Value = '''0.7066'
$ws.Range('E14').Value = '  -2.97%  '
$ws.Range('D15').Value = '''89.57'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').Value = '''29.246.50'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').Value = '''0.000007897'
$ws.Range('E17').Value = '  +2.96%  '
$ws.Range('D18').Value = '''5.790'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('D19').Value = '''13.35'
$ws.Range('E19').Value = '  +1.87%  '
$ws.Range('D20').Value = '''237.10'
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D21').Value = '''1.001'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = '''2.113.42'
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '''7.428'
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('D25').Value = '''162.51'
$ws.Range('D26').Value = '''0.1461'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('D27').Value = '''8.966'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').Value = '''1.957'
$ws.Range('E29').Value = '  -0.80%  '
$ws.Range('D30').Value = '''1.429'
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('E32').Value = '  -2.88%  '
$ws.Range('D33').Value = '''4.015'
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('D34').Value = '''0.05221'
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('D35').Value = '''1.168'
$ws.Range('E35').Value = '  -1.76%  '
$ws.Range('D36').Value = '''0.7075'
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('E37').Value = '  -3.30%  '
$ws.Range('D38').Value = '''2.675'
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('D39').Value = '''0.01859'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').Value = '''2.728'
$ws.Range('E40').Value = '  +1.95%  '
$ws.Range('D41').Value = '''1.141.14'
$ws.Range('E41').Value = '  +6.89%  '
$ws.Range('D42').Value = '''0.9225'
$ws.Range('E42').Value = '  -2.56%  '
$ws.Range('D43').Value = '''0.4282'
$ws.Range('E43').Value = '  -0.56%  '
$ws.Range('D44').Value = '''5.874'
$ws.Range('E44').Value = '  -2.38%  '
$ws.Range('D45').Value = '''70.22'
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').Value = '''1.000'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').Value = '''102.53'
$ws.Range('E47').Value = '  -0.60%  '
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('D49').Value = '''2.010.32'
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('D50').Value = '''9.208'
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('E51').Value = '  -1.40%  '
